# Generate Report for Handoff
# Updates the localization-status workbook so the b.md row reflects that
# its handoff package is now ready (new xliff files were generated for
# both zh-cn and de-de), including the associated status/date/error info
# on the Overview sheet and on each language sheet.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbe4923a2fda080e891f7656ff1b1a97dcbafd53/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab2d659baa0f5a4a69c10ca008c20183ab558df7/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-23 16:40:54"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps this a literal text "False" (matching the
# source file's column, which stores True/False as shared-string text
# rather than a real boolean) instead of Excel auto-coercing it to a
# Boolean cell.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-23 16:40:49"
$wsZhCn.Range("P3").Value = $errorDetail
# Column P (Error Detail) widens to fit the longer text, matching
# the width already used by columns G and J on this sheet.
$wsZhCn.Columns.Item(16).ColumnWidth = $wsZhCn.Columns.Item(7).ColumnWidth

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md entry.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-23 16:40:54"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $wsDeDe.Columns.Item(7).ColumnWidth
